$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray <w:bookmarkStart/bookmarkEnd w:name="_GoBack"/>
#    that currently sits between the "Full" run and the "}" run in
#    the "{subjectFull}" line.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Append a new "{date}" run right after "Date of Submission: ",
#    matching the Times New Roman / 24-half-point run formatting
#    already used on that line.
#
#    Executing Find&Replace against the exact run text first merges
#    the appended text into the same run (so it automatically
#    inherits every rPr property, incl. eastAsia/cs fonts). Then we
#    re-select just "{date}" and nudge its font, which makes the
#    engine split it back out into its own <w:r> while keeping the
#    identical formatting.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Date of Submission: ", $true, $false, $false, $false, $false, $true, 1, $false, "Date of Submission: {date}", 2) | Out-Null

$dateRng = $d.Content
$dateRng.Find.Execute("{date}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$dateRng.Font.Name = "Times New Roman"
$dateRng.Font.Size = 12
$dateRng.Font.NameFarEast = "Times New Roman"
$dateRng.Font.NameBi = "Times New Roman"

# ------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark, now collapsed at the very end
#    of the document (right after the new "{date}" run).
#
#    Placing a zero-length range exactly at the document's last (or
#    second-to-last) character position trips up Bookmarks.Add in
#    this host, so a one-character scratch run is appended first to
#    move the real target position away from that edge, the bookmark
#    is added there, and the scratch character is deleted afterwards
#    -- the bookmark stays put once created.
# ------------------------------------------------------------------
$endPos = $d.Content.End
$scratch = $d.Range($endPos - 1, $endPos - 1)
$scratch.InsertAfter("X")

$target = $d.Range($endPos - 1, $endPos - 1)
$d.Bookmarks.Add("_GoBack", $target)

$endPos2 = $d.Content.End
$scratch2 = $d.Range($endPos2 - 2, $endPos2 - 1)
$scratch2.Text = ""
